$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("E2").Value = 25.29000000000051
$ws.Range("G2").Value = [double]"1.156297280147101e-12"
$ws.Range("H2").Value = [double]"1.789934739379574e-11"
$ws.Range("K2").Value = 37.45344619284136
$ws.Range("L2").Value = "[26.29164391562675, 48.615248470055974]"
$ws.Range("M2").Value = [double]"3.632472100889572e-10"
$ws.Range("N2").Value = [double]"7.264944201779144e-10"
$ws.Range("O2").Value = 1.742184514603348
$ws.Range("P2").Value = "[1.415131825941347, 2.0692372032653497]"
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 58.24813392078907
$ws.Range("T2").Value = "[51.64220429992177, 64.85406354165637]"
$ws.Range("W2").Value = 18.27765765765803
$ws.Range("X2").Value = 16.9612612612616
$ws.Range("Y2").Value = 19.59405405405446

# --- Row 3 updates ---
$ws.Range("B3").Value = 1
$ws.Range("E3").Value = 24.50000000000039
$ws.Range("G3").Value = [double]"5.155431637149377e-12"
$ws.Range("H3").Value = [double]"2.491035009072182e-11"
$ws.Range("I3").Value = 0.05290184476899396
$ws.Range("K3").Value = 38.88571068436831
$ws.Range("L3").Value = "[25.84184413885476, 51.92957722988185]"
$ws.Range("M3").Value = [double]"1.822800710371553e-08"
$ws.Range("N3").Value = [double]"1.822800710371553e-08"
$ws.Range("O3").Value = 1.037763339023655
$ws.Range("P3").Value = "[0.672973801669885, 1.4025528763774258]"
$ws.Range("Q3").Value = [double]"7.050193362445611e-08"
$ws.Range("R3").Value = [double]"7.050193362445611e-08"
$ws.Range("S3").Value = 61.48884656923518
$ws.Range("T3").Value = "[54.54952681178537, 68.42816632668499]"
$ws.Range("W3").Value = 20.45345345345378
$ws.Range("X3").Value = 19.03103103103133
$ws.Range("Y3").Value = 21.87587587587623
